$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.241.53"
$ws.Range("E2").Value = "  +3.81%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.632.52"
$ws.Range("E3").Value = "  +4.06%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.72"
$ws.Range("E5").Value = "  +2.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.30"
$ws.Range("E6").Value = "  +2.23%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.525"
$ws.Range("E8").Value = "  +1.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.631.57"
$ws.Range("E9").Value = "  +4.06%  "

$ws.Range("E10").Value = "  +13.98%  "

$ws.Range("E11").Value = "  +0.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.347"
$ws.Range("E12").Value = "  +2.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.00"
$ws.Range("E13").Value = "  +0.18%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.114.17"
$ws.Range("E14").Value = "  +4.11%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000189"
$ws.Range("E15").Value = "  +10.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.71"
$ws.Range("E16").Value = "  +2.02%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "71.178.59"
$ws.Range("E17").Value = "  +3.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.632.91"
$ws.Range("E18").Value = "  +4.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "381.80"

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.55"
$ws.Range("E20").Value = "  +4.17%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.91"
$ws.Range("E21").Value = "  +4.79%  "

$ws.Range("E22").Value = "  -1.92%  "

$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.43"
$ws.Range("E23").Value = "  +2.21%  "

$ws.Range("B24").Value = "NEARProtocol"
$ws.Range("C24").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.47"
$ws.Range("E24").Value = "  +5.49%  "

$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.93"
$ws.Range("E26").Value = "  +13.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.65"
$ws.Range("E27").Value = "  +6.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.767.37"
$ws.Range("E28").Value = "  +2.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0966"
$ws.Range("E30").Value = "  +7.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "542.26"
$ws.Range("E31").Value = "  +5.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.06"
$ws.Range("E32").Value = "  +3.16%  "

$ws.Range("E33").Value = "  +4.98%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("E34").Value = "  +3.32%  "

$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "166.03"
$ws.Range("E36").Value = "  +1.18%  "

$ws.Range("E37").Value = "  -1.83%  "

$ws.Range("E38").Value = "  +4.31%  "

$ws.Range("E39").Value = "  +6.31%  "

$ws.Range("E40").Value = "  +1.89%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.38"
$ws.Range("E41").Value = "  +4.70%  "

$ws.Range("E42").Value = "  +9.38%  "

$ws.Range("E43").Value = "  +0.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.06"
$ws.Range("E44").Value = "  +4.32%  "

$ws.Range("E45").Value = "  +1.30%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.94"
$ws.Range("E46").Value = "  +2.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "154.48"
$ws.Range("E47").Value = "  +0.69%  "

$ws.Range("E48").Value = "  +1.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.535"
$ws.Range("E50").Value = "  +2.53%  "

$ws.Range("E51").Value = "  +1.67%  "
